$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.524.73"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.847.67"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'258.36"
$ws.Range("E5").Value = "  -7.19%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.5264"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("D8").Value = "'0.3290"
$ws.Range("E8").Value = "  -5.70%  "
$ws.Range("D9").Value = "'0.06745"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").Value = "'0.7774"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "'0.07647"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "1.832.43"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'88.75"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "'5.061"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'14.16"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.000007913"
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("D20").Value = "26.588.02"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "2.067.76"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "'4.606"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").Value = "'9.726"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").Value = "'6.005"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("D25").Value = "'2.353"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "'144.67"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'1.639"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").Value = "'17.02"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "'111.14"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").Value = "'4.235"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").Value = "'4.197"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("D32").Value = "'0.08788"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("D33").Value = "'0.04857"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'1.143"
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").Value = "'2.861"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'0.7093"
$ws.Range("E36").Value = "  -2.65%  "
$ws.Range("D37").Value = "'3.117"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("D38").Value = "'0.01812"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'2.228"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").Value = "'0.4965"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").Value = "'113.61"
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("D42").Value = "'0.9042"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "'6.075"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D44").Value = "'7.830"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "'0.9999"
$ws.Range("D46").Value = "'0.4297"
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("D47").Value = "'0.1293"
$ws.Range("E47").Value = "  -4.90%  "
$ws.Range("D48").Value = "'9.214"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "'0.05930"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'35.36"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").Value = "'1.437"
$ws.Range("E51").Value = "  -3.14%  "
